# Sync automático del tracker - mark the Nantes vs PSG prediction as
# Completed and append its final result to the Results sheet.

$wb = $excel.ActiveWorkbook

# --- 1. Predictions sheet: flip status from "Pending" to "Completed" ---
$wsPred = $wb.Worksheets.Item("Predictions")
$wsPred.Range("I29").Value = "Completed"

# --- 2. Results sheet: append the finished result for that match (row 31) ---
$wsRes = $wb.Worksheets.Item("Results")

# Columns A and J hold dates formatted as plain text (e.g. "2025-08-17") in
# this sheet, not real Excel dates. Force text entry (like typing a leading
# apostrophe in Excel) so the value isn't auto-converted into a date serial,
# then drop back to the Normal style so no extra formatting is left behind.
$wsRes.Range("A31").NumberFormat = "@"
$wsRes.Range("A31").Value = "2025-08-17"
$wsRes.Range("A31").Style = "Normal"

$wsRes.Range("B31").Value = "Ligue 1"
$wsRes.Range("C31").Value = "nantes"
$wsRes.Range("D31").Value = "paris saint germain"
$wsRes.Range("E31").Value = "Away Win"
$wsRes.Range("F31").Value = "Draw"
$wsRes.Range("G31").Value = $false
$wsRes.Range("H31").Value = -1
$wsRes.Range("I31").Value = -100

$wsRes.Range("J31").NumberFormat = "@"
$wsRes.Range("J31").Value = "2025-08-17"
$wsRes.Range("J31").Style = "Normal"
